$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.2862826666666667
$ws.Range("H2").Value = 0.8588480000000001
$ws.Range("I2").Value = 0.560705294934871
$ws.Range("J2").Value = 0.560705294934871
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.001396333333333333
$ws.Range("N2").Value = 0.004189
$ws.Range("O2").Value = 0.003932092785750223
$ws.Range("P2").Value = 0.003932092785750224
$ws.Range("Q2").Value = 0.0003997460302222223
$ws.Range("R2").Value = 0.003597714272
$ws.Range("S2").Value = 0.002204745245145357
$ws.Range("T2").Value = 0.002204745245145358

$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.2862826666666667
$ws.Range("H3").Value = 0.8588480000000001
$ws.Range("I3").Value = 0.560705294934871
$ws.Range("J3").Value = 0.560705294934871
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.277895
$ws.Range("N3").Value = 0.833685
$ws.Range("O3").Value = 0.7825559260177072
$ws.Range("P3").Value = 0.7825559260177072
$ws.Range("Q3").Value = 0.07955652165333334
$ws.Range("R3").Value = 0.7160086948800001
$ws.Range("S3").Value = 0.4387832513007896
$ws.Range("T3").Value = 0.4387832513007896

$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.2862826666666667
$ws.Range("H4").Value = 0.8588480000000001
$ws.Range("I4").Value = 0.560705294934871
$ws.Range("J4").Value = 0.560705294934871
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.06612433333333334
$ws.Range("N4").Value = 0.198373
$ws.Range("O4").Value = 0.1862069807084338
$ws.Range("P4").Value = 0.1862069807084338
$ws.Range("Q4").Value = 0.01893025047822222
$ws.Range("R4").Value = 0.170372254304
$ws.Range("S4").Value = 0.1044072400370542
$ws.Range("T4").Value = 0.1044072400370542

$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.2862826666666667
$ws.Range("H5").Value = 0.8588480000000001
$ws.Range("I5").Value = 0.560705294934871
$ws.Range("J5").Value = 0.560705294934871
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.009696333333333333
$ws.Range("N5").Value = 0.029089
$ws.Range("O5").Value = 0.02730500048810892
$ws.Range("P5").Value = 0.02730500048810892
$ws.Range("Q5").Value = 0.002775892163555556
$ws.Range("R5").Value = 0.024983029472
$ws.Range("S5").Value = 0.01531005835188191
$ws.Range("T5").Value = 0.01531005835188191

$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.1092446666666667
$ws.Range("H6").Value = 0.327734
$ws.Range("I6").Value = 0.2139635757784672
$ws.Range("J6").Value = 0.2139635757784672
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.001396333333333333
$ws.Range("N6").Value = 0.004189
$ws.Range("O6").Value = 0.003932092785750223
$ws.Range("P6").Value = 0.003932092785750224
$ws.Range("Q6").Value = 0.0001525419695555556
$ws.Range("R6").Value = 0.001372877726
$ws.Range("S6").Value = 0.0008413246327318322
$ws.Range("T6").Value = 0.0008413246327318324

$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.1092446666666667
$ws.Range("H7").Value = 0.327734
$ws.Range("I7").Value = 0.2139635757784672
$ws.Range("J7").Value = 0.2139635757784672
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.277895
$ws.Range("N7").Value = 0.833685
$ws.Range("O7").Value = 0.7825559260177072
$ws.Range("P7").Value = 0.7825559260177072
$ws.Range("Q7").Value = 0.03035854664333333
$ws.Range("R7").Value = 0.27322691979
$ws.Range("S7").Value = 0.1674384641773783
$ws.Range("T7").Value = 0.1674384641773783

$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.3333333333333333
$ws.Range("G8").Value = 0.1092446666666667
$ws.Range("H8").Value = 0.327734
$ws.Range("I8").Value = 0.2139635757784672
$ws.Range("J8").Value = 0.2139635757784672
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.06612433333333334
$ws.Range("N8").Value = 0.198373
$ws.Range("O8").Value = 0.1862069807084338
$ws.Range("P8").Value = 0.1862069807084338
$ws.Range("Q8").Value = 0.007223730753555557
$ws.Range("R8").Value = 0.06501357678200001
$ws.Range("S8").Value = 0.03984151142728855
$ws.Range("T8").Value = 0.03984151142728856

$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.3333333333333333
$ws.Range("G9").Value = 0.1092446666666667
$ws.Range("H9").Value = 0.327734
$ws.Range("I9").Value = 0.2139635757784672
$ws.Range("J9").Value = 0.2139635757784672
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.009696333333333333
$ws.Range("N9").Value = 0.029089
$ws.Range("O9").Value = 0.02730500048810892
$ws.Range("P9").Value = 0.02730500048810892
$ws.Range("Q9").Value = 0.001059272702888889
$ws.Range("R9").Value = 0.009533454326000002
$ws.Range("S9").Value = 0.005842275541068576
$ws.Range("T9").Value = 0.005842275541068578

$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 0.6666666666666666
$ws.Range("G10").Value = 0.1150486666666666
$ws.Range("H10").Value = 0.345146
$ws.Range("I10").Value = 0.2253311292866618
$ws.Range("J10").Value = 0.2253311292866618
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 0.3333333333333333
$ws.Range("M10").Value = 0.001396333333333333
$ws.Range("N10").Value = 0.004189
$ws.Range("O10").Value = 0.003932092785750223
$ws.Range("P10").Value = 0.003932092785750224
$ws.Range("Q10").Value = 0.0001606462882222222
$ws.Range("R10").Value = 0.001445816594
$ws.Range("S10").Value = 0.0008860229078730338
$ws.Range("T10").Value = 0.000886022907873034

$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 0.6666666666666666
$ws.Range("G11").Value = 0.1150486666666666
$ws.Range("H11").Value = 0.345146
$ws.Range("I11").Value = 0.2253311292866618
$ws.Range("J11").Value = 0.2253311292866618
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 0.277895
$ws.Range("N11").Value = 0.833685
$ws.Range("O11").Value = 0.7825559260177072
$ws.Range("P11").Value = 0.7825559260177072
$ws.Range("Q11").Value = 0.03197144922333333
$ws.Range("R11").Value = 0.28774304301
$ws.Range("S11").Value = 0.1763342105395393
$ws.Range("T11").Value = 0.1763342105395393

$ws.Range("E12").Value = 2
$ws.Range("F12").Value = 0.6666666666666666
$ws.Range("G12").Value = 0.1150486666666666
$ws.Range("H12").Value = 0.345146
$ws.Range("I12").Value = 0.2253311292866618
$ws.Range("J12").Value = 0.2253311292866618
$ws.Range("K12").Value = 2
$ws.Range("L12").Value = 0.6666666666666666
$ws.Range("M12").Value = 0.06612433333333334
$ws.Range("N12").Value = 0.198373
$ws.Range("O12").Value = 0.1862069807084338
$ws.Range("P12").Value = 0.1862069807084338
$ws.Range("Q12").Value = 0.007607516384222222
$ws.Range("R12").Value = 0.06846764745799999
$ws.Range("S12").Value = 0.04195822924409104
$ws.Range("T12").Value = 0.04195822924409104

$ws.Range("E13").Value = 2
$ws.Range("F13").Value = 0.6666666666666666
$ws.Range("G13").Value = 0.1150486666666666
$ws.Range("H13").Value = 0.345146
$ws.Range("I13").Value = 0.2253311292866618
$ws.Range("J13").Value = 0.2253311292866618
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 0.6666666666666666
$ws.Range("M13").Value = 0.009696333333333333
$ws.Range("N13").Value = 0.029089
$ws.Range("O13").Value = 0.02730500048810892
$ws.Range("P13").Value = 0.02730500048810892
$ws.Range("Q13").Value = 0.001115550221555555
$ws.Range("R13").Value = 0.010039951994
$ws.Range("S13").Value = 0.006152666595158434
$ws.Range("T13").Value = 0.006152666595158436
